# Apply the "complete IPL scrapping" edit:
# 1. Rename the sheet from "Sheet1" to "Kamlesh Nagarkoti"
# 2. Insert a new first column "matchNo" (shifting all existing columns right by one)
# 3. Fill the new column with header "matchNo" and value "15th"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Kamlesh Nagarkoti"

# Insert a new column at A, shifting existing data (A:L) to (B:M)
$ws.Columns.Item(1).Insert()

# Populate the new column A with header + value
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "15th"
